$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status draft -> active
$ws.Range("B6").Value = "active"

# Experimental: clear the "true" value (row stays, cell becomes empty)
$ws.Range("B7").Value = ""

# Date 2024-01-18 -> 2025-11-18
# Assigning the literal string via .Value auto-converts date-like text to a
# date serial number (like real Excel typing into a General cell). Route the
# text through a formula + paste-values so it lands as plain text, matching
# the original cell's type/style (no NumberFormat change).
$ws.Range("B8").Formula = '="2025-11-18"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
